$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Analysis_info")
$ws2 = $wb.Worksheets.Item("File_paths")

# --- Swap the "last name" author columns (A <-> B) for header + 4 data rows ---
# Row 5 (header): also fix capitalization "last author" -> "Last author"
$ws1.Range("A5").Value = "Last name (Last author)"
$ws1.Range("B5").Value = "Last name (first author)"

# Rows 6-9 (data): swap the first/last author last-name values between columns A and B
for ($r = 6; $r -le 9; $r++) {
    $colA = $ws1.Cells.Item($r, 1).Value2
    $colB = $ws1.Cells.Item($r, 2).Value2
    $ws1.Cells.Item($r, 1).Value = $colB
    $ws1.Cells.Item($r, 2).Value = $colA
}

# --- Add four new blank (but formatted) rows below the table ---
$ws1.Range("A6").Copy()
$ws1.Range("A16").PasteSpecial(-4122)
$ws1.Range("A17").PasteSpecial(-4122)
$ws1.Range("A18").PasteSpecial(-4122)
$ws1.Range("A19").PasteSpecial(-4122)

# --- Update the active sheet / selection state ---
$ws2.Range("A2").Select()
$ws1.Activate()
$ws1.Range("B11").Select()
